{"js": "// The date field \"\u041a\u04af\u043d\u0456: 08.09.2023\u0436\" must become \"\u041a\u04af\u043d\u0456: 11.01.2024\u0436\".\n// The old date digits are split across many runs (one run per char/pair),\n// so a plain text search across the body (which matches across run\n// boundaries) is used to locate the date and replace it in one shot.\nconst body = context.document.body;\n\nconst results = body.search(\"08.09.2023\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find date '08.09.2023' to replace.\");\n}\n\n// Replace the matched range's text in place, preserving its formatting.\nresults.items[0].insertText(\"11.01.2024\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The date field \"\u041a\u04af\u043d\u0456: 08.09.2023\u0436\" must become \"\u041a\u04af\u043d\u0456: 11.01.2024\u0436\".\n# The old date digits are split across many small runs (one run per\n# character/pair), so Find/Replace (which matches across run boundaries)\n# is used to locate \"08.09.2023\" and replace it with \"11.01.2024\" in a\n# single operation, preserving the surrounding run formatting.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"08.09.2023\"\n$find.Replacement.Text = \"11.01.2024\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
